# Fixed Isuues in afact and dct catch trials
# Update the "Phases" column (column G) values on Sheet1 with corrected labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = "Dichotic_and_AFACT"
    3  = "Digit_before_and_AFACT"
    4  = "dichotic_phase"
    5  = "Dichotic_and_AFACT"
    6  = "MAB_and_AFACT"
    7  = "MAB_phase"
    8  = "dichotic_phase"
    9  = "MAB_phase"
    10 = "dichotic_phase"
    11 = "MAB_and_Digit_after"
    12 = "Dichotic_and_AFACT"
    13 = "dichotic_phase"
    14 = "Dichotic_and_AFACT"
    15 = "Digit_before_and_AFACT"
    16 = "MAB_and_AFACT"
    17 = "Dichotic_and_AFACT"
    18 = "dichotic_phase"
    19 = "dichotic_phase"
    20 = "MAB_and_AFACT"
    21 = "dichotic_phase"
    22 = "Dichotic_and_AFACT"
    23 = "dichotic_phase"
    24 = "dichotic_phase"
    25 = "MAB_phase"
    26 = "MAB_and_Digit_after"
    27 = "MAB_and_Digit_after"
    28 = "MAB_and_Digit_after"
    29 = "MAB_phase"
    30 = "Dichotic_and_AFACT"
    31 = "Dichotic_and_AFACT"
    32 = "dichotic_phase"
    33 = "Digit_before_and_AFACT"
    34 = "MAB_and_Digit_after"
    35 = "MAB_phase"
    36 = "MAB_and_AFACT"
    37 = "Digit_before_and_AFACT"
    38 = "Digit_before_and_AFACT"
    39 = "Dichotic_and_AFACT"
    40 = "Dichotic_and_AFACT"
    41 = "MAB_and_AFACT"
    42 = "Digit_before_and_AFACT"
    43 = "dichotic_phase"
    44 = "dichotic_phase"
    45 = "Dichotic_and_AFACT"
    46 = "Dichotic_and_AFACT"
    47 = "MAB_and_AFACT"
    48 = "Dichotic_and_AFACT"
    49 = "MAB_phase"
    50 = "Dichotic_and_AFACT"
    51 = "Digit_before_and_AFACT"
    52 = "dichotic_phase"
    53 = "MAB_phase"
    54 = "MAB_and_Digit_after"
    55 = "Digit_before_and_AFACT"
    56 = "MAB_phase"
    57 = "MAB_and_AFACT"
    58 = "dichotic_phase"
    59 = "dichotic_phase"
    60 = "Digit_before_and_AFACT"
    61 = "Dichotic_and_AFACT"
    62 = "MAB_phase"
    63 = "MAB_and_Digit_after"
    64 = "Digit_before_and_AFACT"
    65 = "MAB_and_AFACT"
    66 = "dichotic_phase"
    67 = "Dichotic_and_AFACT"
    68 = "MAB_and_Digit_after"
    69 = "dichotic_phase"
    70 = "Dichotic_and_AFACT"
    71 = "MAB_phase"
    72 = "MAB_and_Digit_after"
    73 = "MAB_and_AFACT"
    74 = "MAB_and_Digit_after"
    75 = "dichotic_phase"
    76 = "dichotic_phase"
    77 = "Dichotic_and_AFACT"
    78 = "MAB_and_AFACT"
    79 = "Dichotic_and_AFACT"
    80 = "dichotic_phase"
    81 = "Dichotic_and_AFACT"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
